$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Rename the table column header "23/out" -> "Coluna1" (this is the
# header cell B2, part of Tabela2). The engine does not propagate a
# ListColumn / header-cell rename into existing structured-reference
# formulas automatically, so the SUBTOTAL formula below is rewritten
# to match by hand (mirrors what real Excel does under the hood).
$ws.Range("B2").Value = "Coluna1"
$ws.Range("B9").Formula = "=SUBTOTAL(109,Tabela2[Coluna1])"

# Clear the price values that were entered under that column.
$ws.Range("B3:B8").ClearContents()

# Update the active selection to C6.
$ws.Range("C6").Select()

$wb.Application.Calculate()
